$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column) on the "Repayment schedule" sheet,
# shifting the old N/O/P columns (Late / heading / Outstanding) to O/P/Q.
$ws.Columns("N:N").Insert()

# The freshly inserted column keeps the width used by the "Outstanding" block
# (approx. 10.71 characters in the original workbook).
$ws.Columns("N:N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab and update its selection.
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
